# Refresh Market Board valuation columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) for a handful of leves across each crafting-job worksheet,
# per the latest scheduled price pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 38909.816
$ws.Range("J19").Value = 53263.625
$ws.Range("L19").Value = 53263.625
$ws.Range("N19").Value = -53613.625
$ws.Range("H33").Value = 22729864
$ws.Range("I33").Value = 35714644
$ws.Range("K33").Value = 35714644
$ws.Range("M33").Value = -35714415
$ws.Range("N33").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H62").Value = 3999.3333
$ws.Range("J62").Value = 3999.3333
$ws.Range("L62").Value = 3999.3333
$ws.Range("N62").Value = -5247.3333
$ws.Range("H65").Value = 3999.3333
$ws.Range("J65").Value = 3999.3333
$ws.Range("L65").Value = 19996.6665
$ws.Range("N65").Value = -26236.6665
$ws.Range("H92").Value = 1927.1765
$ws.Range("I92").Value = 2041.0625
$ws.Range("K92").Value = 2041.0625
$ws.Range("M92").Value = -793.0625
$ws.Range("N92").ClearContents()
$ws.Range("H129").Value = 1699.6666
$ws.Range("J129").Value = 3249.5
$ws.Range("L129").Value = 9748.5
$ws.Range("N129").Value = -19748.5
$ws.Range("H141").Value = 4072.3635
$ws.Range("I141").Value = 4072.3635
$ws.Range("K141").Value = 12217.0905
$ws.Range("M141").Value = -7037.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1857.5
$ws.Range("I2").Value = 1715.25
$ws.Range("J2").Value = 1999.75
$ws.Range("K2").Value = 1715.25
$ws.Range("L2").Value = 1999.75
$ws.Range("M2").Value = -1602.25
$ws.Range("N2").Value = -2225.75
$ws.Range("H116").Value = 1857.5
$ws.Range("I116").Value = 1715.25
$ws.Range("J116").Value = 1999.75
$ws.Range("K116").Value = 1715.25
$ws.Range("L116").Value = 1999.75
$ws.Range("M116").Value = 578.75
$ws.Range("N116").Value = -6587.75
$ws.Range("H122").Value = 1813.3334
$ws.Range("I122").Value = 1470
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4410
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1960
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1857.5
$ws.Range("I3").Value = 1715.25
$ws.Range("J3").Value = 1999.75
$ws.Range("K3").Value = 1715.25
$ws.Range("L3").Value = 1999.75
$ws.Range("M3").Value = -1601.25
$ws.Range("N3").Value = -2227.75
$ws.Range("H6").Value = 44909.25
$ws.Range("I6").Value = 26647
$ws.Range("K6").Value = 26647
$ws.Range("M6").Value = -26534
$ws.Range("N6").ClearContents()
$ws.Range("H86").Value = 1260
$ws.Range("I86").Value = 1260
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1260
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -137
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1260
$ws.Range("I89").Value = 1260
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6300
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -684
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 1338.4117
$ws.Range("I107").Value = 1480
$ws.Range("K107").Value = 1480
$ws.Range("M107").Value = 440
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 67503.39999999999
$ws.Range("I134").Value = 109002.6
$ws.Range("K134").Value = 327007.8
$ws.Range("M134").Value = -324472.8
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 754.65
$ws.Range("I7").Value = 916.1539
$ws.Range("J7").Value = 454.7143
$ws.Range("K7").Value = 916.1539
$ws.Range("L7").Value = 454.7143
$ws.Range("M7").Value = -803.1539
$ws.Range("N7").Value = -680.7143
$ws.Range("H94").Value = 997.3333
$ws.Range("J94").Value = 997.6
$ws.Range("L94").Value = 997.6
$ws.Range("N94").Value = -1899.6
$ws.Range("H105").Value = 19668.334
$ws.Range("I105").Value = 52005
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 52005
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -50258
$ws.Range("N105").Value = -6994
$ws.Range("H122").Value = 3718.0908
$ws.Range("I122").Value = 2495
$ws.Range("J122").Value = 3840.4
$ws.Range("K122").Value = 7485
$ws.Range("L122").Value = 11521.2
$ws.Range("M122").Value = -5035
$ws.Range("N122").Value = -16421.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 118
$ws.Range("I7").Value = 118
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 354
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -242
$ws.Range("N7").ClearContents()
$ws.Range("H63").Value = 3036.6667
$ws.Range("I63").Value = 3036.6667
$ws.Range("K63").Value = 9110.000100000001
$ws.Range("M63").Value = -8361.000100000001
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 3036.6667
$ws.Range("I66").Value = 3036.6667
$ws.Range("K66").Value = 27330.0003
$ws.Range("M66").Value = -23586.0003
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("H103").Value = 1525
$ws.Range("J103").Value = 1445.8334
$ws.Range("L103").Value = 4337.5002
$ws.Range("N103").Value = -6095.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 424
$ws.Range("I107").Value = 424
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 424
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1496
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 5150
$ws.Range("I113").Value = 4833.3335
$ws.Range("K113").Value = 4833.3335
$ws.Range("M113").Value = -2663.3335
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2142.9565
$ws.Range("I16").Value = 2333.1333
$ws.Range("J16").Value = 1786.375
$ws.Range("K16").Value = 2333.1333
$ws.Range("L16").Value = 1786.375
$ws.Range("M16").Value = -2163.1333
$ws.Range("N16").Value = -2126.375
$ws.Range("H42").Value = 31018.666
$ws.Range("I42").Value = 20000
$ws.Range("J42").Value = 36528
$ws.Range("K42").Value = 20000
$ws.Range("L42").Value = 36528
$ws.Range("M42").Value = -19437
$ws.Range("N42").Value = -37654
$ws.Range("H49").Value = 31018.666
$ws.Range("I49").Value = 20000
$ws.Range("J49").Value = 36528
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 36528
$ws.Range("M49").Value = -19853
$ws.Range("N49").Value = -36822
$ws.Range("H93").Value = 6954.375
$ws.Range("I93").Value = 4376.857
$ws.Range("J93").Value = 24997
$ws.Range("K93").Value = 4376.857
$ws.Range("L93").Value = 24997
$ws.Range("M93").Value = -3128.857
$ws.Range("N93").Value = -27493
$ws.Range("H132").Value = 2503055
$ws.Range("I132").Value = 2856.4614
$ws.Range("J132").Value = 4669894
$ws.Range("K132").Value = 8569.3842
$ws.Range("L132").Value = 14009682
$ws.Range("M132").Value = -6039.3842
$ws.Range("N132").Value = -14014742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 755.2
$ws.Range("I81").Value = 755.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1510.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -449.4000000000001
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 755.2
$ws.Range("I84").Value = 755.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7552
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2248
$ws.Range("N84").ClearContents()
$ws.Range("H99").Value = 47886.4
$ws.Range("I99").Value = 44716
$ws.Range("J99").Value = 50000
$ws.Range("K99").Value = 44716
$ws.Range("L99").Value = 50000
$ws.Range("M99").Value = -41721
$ws.Range("N99").Value = -55990
$ws.Range("H122").Value = 3320.075
$ws.Range("I122").Value = 1943.3214
$ws.Range("K122").Value = 5829.9642
$ws.Range("M122").Value = -3379.9642
$ws.Range("N122").ClearContents()
